# Update "want to go" counts (column F) across sheets, matching the
# upstream generated-data refresh (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 125
$ws1.Range("F3").Value = 1289
$ws1.Range("F4").Value = 934
$ws1.Range("F5").Value = 975
$ws1.Range("F6").Value = 1739
$ws1.Range("F8").Value = 1160
$ws1.Range("F13").Value = 48
$ws1.Range("F15").Value = 646
$ws1.Range("F16").Value = 140
$ws1.Range("F27").Value = 849
$ws1.Range("F28").Value = 301

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 299

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 299
$ws4.Range("F3").Value = 125
$ws4.Range("F4").Value = 1289
$ws4.Range("F5").Value = 934
$ws4.Range("F6").Value = 975
$ws4.Range("F7").Value = 1739
$ws4.Range("F9").Value = 1160
$ws4.Range("F15").Value = 48
$ws4.Range("F17").Value = 646
$ws4.Range("F18").Value = 140
$ws4.Range("F35").Value = 849
$ws4.Range("F36").Value = 301
